$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# Name: CompetencemetierVs -> CompetencemetiereVs
$meta.Range("B4").Value = "CompetencemetiereVs"

# Date: 2025-07-21T14:08:48+00:00 -> 2025-07-22T15:23:11+00:00
$meta.Range("B8").Value = "2025-07-22T15:23:11+00:00"

# --- Include #0 sheet ---
$inc0 = $wb.Worksheets.Item("Include #0")
$inc0.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R01-EnsembleSavoirFaire-CISIS"

# --- Include #1 sheet ---
$inc1 = $wb.Worksheets.Item("Include #1")
$inc1.Range("B4").Value = "https://smt.esante.gouv.fr/fhir/CodeSystem/tre-r394-competence-metier"
